$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite vendor names in column A (rows 2-8) to the dashed/uppercase form
$ws.Range("A2").Value = "VENDOR-01"
$ws.Range("A3").Value = "VENDOR-02"
$ws.Range("A4").Value = "VENDOR-03"
$ws.Range("A5").Value = "VENDOR-04"
$ws.Range("A6").Value = "VENDOR-05"
$ws.Range("A7").Value = "VENDOR-06"
$ws.Range("A8").Value = "VENDOR-07"

# Update the selection to match the new active cell / selected range
$ws.Range("A2:A8").Select()
$excel.ActiveCell = $ws.Range("A2")
